$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Cell values — written in the exact order new text needs to appear so the
#    rebuilt shared-string table lines up with the target ("Elevation",
#    "Fire" and "P" already exist in A3 / A4 / F2 and are left untouched so
#    they keep their original table slots).
# ---------------------------------------------------------------------------

# new strings, first-use order
$ws.Range("C2").Value2 = "Chisq"
$ws.Range("A5").Value2 = "Elevation:Fire"
$ws.Range("E1").Value2 = "Tree Height (m)"
$ws.Range("G1").Value2 = "Canopy Spread (m)"
$ws.Range("I1").Value2 = "DBH (cm)"
$ws.Range("K1").Value2 = "Distance Between Neighbors (m)"
$ws.Range("B2").Value2 = "df"
$ws.Range("D2").Value2 = "Pr(>Chisq)"
$ws.Range("C1").Value2 = "Slope"

# repeated references to strings already introduced above
$ws.Range("E2").Value2 = "Chisq"
$ws.Range("G2").Value2 = "Chisq"
$ws.Range("I2").Value2 = "Chisq"
$ws.Range("K2").Value2 = "Chisq"
$ws.Range("H2").Value2 = "P"
$ws.Range("J2").Value2 = "P"
$ws.Range("L2").Value2 = "P"

# row labels (A3 "Elevation" / A4 "Fire" already correct and left as-is)
$ws.Range("B3").Value2 = 1
$ws.Range("B4").Value2 = 1
$ws.Range("B5").Value2 = 1

# numeric statistics grid
$ws.Range("C3").Value2 = 1.47821059364559004301
$ws.Range("D3").Value2 = 0.22405470151761300301
$ws.Range("E3").Value2 = 3.45065688668170000852
$ws.Range("F3").Value2 = 0.06322663698520550324
$ws.Range("G3").Value2 = 7.94786062419213035213
$ws.Range("H3").Value2 = 0.00481442558796878022
$ws.Range("I3").Value2 = 3.43340370406360984035
$ws.Range("J3").Value2 = 0.06389029910355209507
$ws.Range("K3").Value2 = 0.68329848469389597998
$ws.Range("L3").Value2 = 0.40845321160451397713

$ws.Range("C4").Value2 = 1.54169292873014995848
$ws.Range("D4").Value2 = 0.21436599192908800671
$ws.Range("E4").Value2 = 0.09724470921803110568
$ws.Range("F4").Value2 = 0.75516148506394997764
$ws.Range("G4").Value2 = 0.01186947664114829971
$ws.Range("H4").Value2 = 0.91324442655978599337
$ws.Range("I4").Value2 = 0.1569679468529419919
$ws.Range("J4").Value2 = 0.69196353593496096046
$ws.Range("K4").Value2 = 0.01240324767453659938
$ws.Range("L4").Value2 = 0.911323052614901985

$ws.Range("C5").Value2 = 0.25986113488504197777
$ws.Range("D5").Value2 = 0.61021557302515705423
$ws.Range("E5").Value2 = 6.59332984604079985758
$ws.Range("F5").Value2 = 0.01023615364476869981
$ws.Range("G5").Value2 = 0.06835770668363840008
$ws.Range("H5").Value2 = 0.79374313748962599391
$ws.Range("I5").Value2 = 0.06056643781997580328
$ws.Range("J5").Value2 = 0.80560289209100099939
$ws.Range("K5").Value2 = 2.92963093545583008748
$ws.Range("L5").Value2 = 0.08696737135663129559

# clear the stray old row-6 "Residuals"/36 leftovers before restyling row 6/7
$ws.Range("A6").ClearContents() | Out-Null
$ws.Range("B6").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# 2) Number formats — "0.000" (numFmtId 164) on every results column.
# ---------------------------------------------------------------------------
$ws.Range("C3:L5").NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# 3) Fonts — label/header cells keep the small (10pt) font; data cells move
#    to the workbook's default (12pt) font.
# ---------------------------------------------------------------------------
$ws.Range("C3:L5").Font.Size = 12
$ws.Range("K6:L6").Font.Size = 12
$ws.Range("A7:I7").Font.Size = 12
$ws.Range("K7:L7").Font.Size = 12

# ---------------------------------------------------------------------------
# 4) Borders — the old boxed border (style index 10) is no longer used
#    anywhere, so drop it from every cell that had it.
# ---------------------------------------------------------------------------
$ws.Range("A1:L1").Borders.LineStyle = -4142
$ws.Range("A2:B2").Borders.LineStyle = -4142
$ws.Range("E2:L2").Borders.LineStyle = -4142
$ws.Range("A3:B5").Borders.LineStyle = -4142
$ws.Range("C3:L5").Borders.LineStyle = -4142
$ws.Range("A6:J6").Borders.LineStyle = -4142
$ws.Range("J7").Borders.LineStyle = -4142
$ws.Range("K6:L6").Borders.LineStyle = -4142
$ws.Range("A7:I7").Borders.LineStyle = -4142
$ws.Range("K7:L7").Borders.LineStyle = -4142

# ---------------------------------------------------------------------------
# 5) Column widths (best-effort match to the saved "best fit" values).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 11.498697916666666
$ws.Columns.Item(2).ColumnWidth = 1.9986979166666667
$ws.Columns.Item(3).ColumnWidth = 4.830729166666667
$ws.Columns.Item(4).ColumnWidth = 8.666666666666666
$ws.Columns.Item(5).ColumnWidth = 13.330729166666666
$ws.Columns.Item(6).ColumnWidth = 11.330729166666666
$ws.Columns.Item(7).ColumnWidth = 15.998697916666666
$ws.Columns.Item(8).ColumnWidth = 11.330729166666666
$ws.Columns.Item(9).ColumnWidth = 11.330729166666666
$ws.Columns.Item(10).ColumnWidth = 11.330729166666666
$ws.Columns.Item(11).ColumnWidth = 27.998697916666668
$ws.Columns.Item(12).ColumnWidth = 11.330729166666666

# ---------------------------------------------------------------------------
# 6) Selection
# ---------------------------------------------------------------------------
$ws.Range("C3:D5").Select() | Out-Null
